$d = $word.ActiveDocument

# --- Fill in the (previously empty) first paragraph with the project idea text ---
$p = $d.Paragraphs(1)
$r = $p.Range
$r.InsertAfter("Educational program where you could select a subject and that would bring up various questions for that subject")
$r.Collapse(0)
$r.InsertAfter(". It would provide appropriate responses for right and wrong answers. We could either import the questions from ")
$r.Collapse(0)
$r.InsertAfter("a module or just use a few different questions for each subject, whichever is most effective.")

# --- Apply the "List Paragraph" style to the paragraph ---
$p.Style = "List Paragraph"

# Match the style's generated formatting to a first-use built-in "List Paragraph" style
$style = $d.Styles("List Paragraph")
$style.Priority = 34
$style.NoSpaceBetweenParagraphsOfSameStyle = $true
$style.ParagraphFormat.LeftIndent = 36

# --- Turn the paragraph into item 1 of a new numbered list (numId 1) ---
$p.Range.ListFormat.ApplyNumberDefault()

# Shape the generated list definition to match the classic "1)  a.  i." numbering
# gallery entry (decimal/lowerLetter/lowerRoman repeating every three levels).
$lt = $p.Range.ListFormat.ListTemplate
$formats = @("%1)", "%2.", "%3.", "%4.", "%5.", "%6.", "%7.", "%8.", "%9.")
$styles  = @(0, 4, 2, 0, 4, 2, 0, 4, 2)
for ($i = 1; $i -le $lt.ListLevels.Count; $i++) {
    $lvl = $lt.ListLevels.Item($i)
    $lvl.NumberStyle = $styles[$i - 1]
    $lvl.NumberFormat = $formats[$i - 1]
}
